# NATMI LR-pair recompute with refreshed TPM input (Ren1 -> Atp6ap2).
# Rows 2-5 (sending cluster = ECs) get refreshed stats, and the same
# four target clusters are appended as rows 6-9 for sending cluster =
# MuSCs, which NATMI now also reports. Used range grows to A1:T9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns A..T in order.
$rowsData = @(
    ,@(2, "ECs", "Ren1", "Atp6ap2", "ECs", 3, 1, 1.307223, 3.921669, 0.8940139324750501, 0.89401393247505, 3, 1, 8.382792333333333, 25.148377, 0.0868480383930262, 0.08684803839302623, 10.958178942357, 98.623610481213, 0.07764335633149348, 0.0776433563314935)
    ,@(3, "ECs", "Ren1", "Atp6ap2", "FAPs", 3, 1, 1.307223, 3.921669, 0.8940139324750501, 0.89401393247505, 3, 1, 20.20471733333333, 60.614152, 0.2093264388416289, 0.2093264388416289, 26.412071206632, 237.708640859688, 0.1871407527598027, 0.1871407527598027)
    ,@(4, "ECs", "Ren1", "Atp6ap2", "MuSCs", 3, 1, 1.307223, 3.921669, 0.8940139324750501, 0.89401393247505, 3, 1, 16.00744233333333, 48.022327, 0.1658415133119111, 0.1658415133119112, 20.925296789307, 188.327671103763, 0.148264623483595, 0.148264623483595)
    ,@(5, "ECs", "Ren1", "Atp6ap2", "Resolving-Mac", 3, 1, 1.307223, 3.921669, 0.8940139324750501, 0.89401393247505, 3, 1, 51.92757733333334, 155.782732, 0.5379840094534337, 0.5379840094534337, 67.880923424412, 610.928310819708, 0.4809651999001588, 0.4809651999001587)
    ,@(6, "MuSCs", "Ren1", "Atp6ap2", "ECs", 1, 0.3333333333333333, 0.1549723333333334, 0.464917, 0.1059860675249499, 0.1059860675249499, 3, 1, 8.382792333333333, 25.148377, 0.0868480383930262, 0.08684803839302623, 1.299100887745444, 11.691907989709, 0.00920468206153272, 0.00920468206153272)
    ,@(7, "MuSCs", "Ren1", "Atp6ap2", "FAPs", 1, 0.3333333333333333, 0.1549723333333334, 0.464917, 0.1059860675249499, 0.1059860675249499, 3, 1, 20.20471733333333, 60.614152, 0.2093264388416289, 0.2093264388416289, 3.131172189487112, 28.180549705384, 0.02218568608182618, 0.02218568608182618)
    ,@(8, "MuSCs", "Ren1", "Atp6ap2", "MuSCs", 1, 0.3333333333333333, 0.1549723333333334, 0.464917, 0.1059860675249499, 0.1059860675249499, 3, 1, 16.00744233333333, 48.022327, 0.1658415133119111, 0.1658415133119112, 2.480710689095445, 22.326396201859, 0.01757688982831609, 0.0175768898283161)
    ,@(9, "MuSCs", "Ren1", "Atp6ap2", "Resolving-Mac", 1, 0.3333333333333333, 0.1549723333333334, 0.464917, 0.1059860675249499, 0.1059860675249499, 3, 1, 51.92757733333334, 155.782732, 0.5379840094534337, 0.5379840094534337, 8.047337823693779, 72.42604041324401, 0.05701880955327494, 0.05701880955327492)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $rowsData) {
    $r = $entry[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "$($columns[$i])$r"
        $ws.Range($cellRef).Value = $entry[$i + 1]
    }
}
